# RPA datasets push 2024-03-23
# Shift rows 5..19 down to 7..21 (copy whole rows, bottom-up to avoid clobbering),
# then write the new top rows 3..6 with fresh data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Shift existing rows 5-19 down by two rows into 7-21, processing from the
#    bottom up so the source rows are not overwritten before they are read.
for ($oldRow = 19; $oldRow -ge 5; $oldRow--) {
    $newRow = $oldRow + 2
    $src = $ws.Range("A" + $oldRow + ":F" + $oldRow)
    $dst = $ws.Range("A" + $newRow + ":F" + $newRow)
    $src.Copy($dst)
}

# 2) Row 2 (노브랜드) is unchanged.

# 3) Row 3: new entry - 아이씨티케이
$ws.Range("A3").Value2 = "아이씨티케이"
$ws.Range("B3").Value2 = "2024.04.24~04.30"
$ws.Range("C3").Value2 = "13,000~16,000"
$ws.Range("D3").Value2 = $ws.Range("D2").Value2
$ws.Range("E3").Value2 = 25610
$ws.Range("F3").Value2 = "NH투자증권"

# 4) Row 4: 이노그리드 (was row 3)
$ws.Range("A4").Value2 = "이노그리드"
$ws.Range("B4").Value2 = "2024.04.18~04.24"
$ws.Range("C4").Value2 = "29,000~35,000"
$ws.Range("D4").Value2 = $ws.Range("D2").Value2
$ws.Range("E4").Value2 = 17400
$ws.Range("F4").Value2 = "한국투자증권"

# 5) Row 5: 코칩 (was row 4)
$ws.Range("A5").Value2 = "코칩"
$ws.Range("B5").Value2 = "2024.04.15~04.19"
$ws.Range("C5").Value2 = "11,000~14,000"
$ws.Range("D5").Value2 = $ws.Range("D2").Value2
$ws.Range("E5").Value2 = 16500
$ws.Range("F5").Value2 = "한국투자증권"

# 6) Row 6: new entry - 유안타스팩16호
$ws.Range("A6").Value2 = "유안타스팩16호"
$ws.Range("B6").Value2 = "2024.04.15~04.16"
$ws.Range("C6").Value2 = "2,000~2,000"
$ws.Range("D6").Value2 = $ws.Range("D2").Value2
$ws.Range("E6").Value2 = 10300
$ws.Range("F6").Value2 = "유안타증권"
